# adicionando dados do campeonato mundial de kite 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fórmula Kite Masc. (row 11) and Fórmula Kite Fem. (row 12) now have a
# result for the 2024 World Championship (column C).
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1

# An underline style was applied to an empty cell below the table (D17),
# and that cell became the active selection.
$ws.Range("D17").Font.Underline = $true
$ws.Range("D17").Select()
